$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = ""
$ws.Range("C1").Value = "Nr of points"
$ws.Range("D1").Value = "Points"
$ws.Range("E1").Value = "%-age"
$ws.Range("F1").Value = "Column 6"
$ws.Range("G1").Value = "Column 7"
$ws.Range("H1").Value = "Column 8"
$ws.Range("I1").Value = "Column 9"
$ws.Range("J1").Value = "Column 10"
